$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update existing row 4 values (free-float calculation adjustment) ---
$ws.Range("A4").Value = 44298
$ws.Range("B4").Value = 0.2655415777506262
$ws.Range("H4").Value = 0.2921338026027105

# --- 2) Add a new row of data in row 5, reusing row 4's formatting ---
$ws.Range("A4:J4").Copy()
$ws.Range("A5:J5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A5").Value = 44298
$ws.Range("B5").Value = 0.28459279237716756
$ws.Range("C5").Value = 0.9965552505516981
$ws.Range("D5").Value = 0.007111272306325617
$ws.Range("E5").Value = 0.020880696419320013
$ws.Range("F5").Value = 0.6445328055569338
$ws.Range("G5").Value = 0.2974762569832402
$ws.Range("H5").Value = 0.2921338026027105
$ws.Range("I5").Value = 0.2468135162241888
$ws.Range("J5").Value = 0.3590013491833052

# --- 3) Add a brand new row 6 with contribution data, reusing row 4's formatting ---
$ws.Range("A4:J4").Copy()
$ws.Range("A6:J6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A6").Value = 44306
$ws.Range("B6").Value = 0.28459279237716756
$ws.Range("C6").Value = 0.9965552505516981
$ws.Range("D6").Value = 0.007111272306325617
$ws.Range("E6").Value = 0.020880696419320013
$ws.Range("F6").Value = 0.6445328055569338
$ws.Range("G6").Value = 0.29749133191451627
$ws.Range("H6").Value = 0.32864406481746694
$ws.Range("I6").Value = 0.2468135162241888
$ws.Range("J6").Value = 0.3590013491833052

# --- 4) Adjust the contribution base table: restyle L11 (bold, 0.000% format) ---
$ws.Range("L11").NumberFormat = "0.000%"
$ws.Range("L11").Font.Bold = $true

# --- 5) Update the active selection to reflect where the editor left off ---
[void]$ws.Range("I13").Select()
